# Adds 4 new data rows (31-34) to the active worksheet, matching the target diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @{ 'A' = 'run_2025-09-29T21-42-48.129862p00-00_96d301d8.json'; 'B' = '96d301d8'; 'C' = 'success'; 'D' = '2025-09-29T21-42-48.129862p00-00'; 'E' = 1; 'F' = 'test'; 'G' = 'openai/gpt-5-mini'; 'H' = 'realpariwise-alltools-gpt5mini-1-test'; 'I' = 'coherence_check, counterfactual_pairs'; 'J' = 0.6666666666666666; 'K' = 0.3333333333333333; 'L' = 0.4444444444444444; 'M' = 0; 'N' = 1; 'O' = 0.3333333333333333; 'P' = 0.5; 'Q' = 3; 'R' = 0.5; 'S' = 0.3333333333333333; 'T' = 0.4; 'U' = 3; 'V' = 0.9166666666666666; 'W' = 0.9777777777777777; 'X' = 0.946236559139785; 'Y' = 45 },
    @{ 'A' = 'run_2025-09-29T22-42-49.468928p00-00_fc95b218.json'; 'B' = 'fc95b218'; 'C' = 'success'; 'D' = '2025-09-29T22-42-49.468928p00-00'; 'E' = 1; 'F' = 'test'; 'G' = 'openai/gpt-5-mini'; 'H' = 'realpariwise-alltools-gpt5mini-1-test'; 'I' = 'coherence_check, counterfactual_pairs'; 'J' = 0; 'K' = 0; 'L' = 0; 'M' = 0; 'N' = 0; 'O' = 0; 'P' = 0; 'Q' = 3; 'R' = 0; 'S' = 0; 'T' = 0; 'U' = 3; 'V' = 0.8823529411764706; 'W' = 1; 'X' = 0.9375; 'Y' = 45 },
    @{ 'A' = 'run_2025-09-29T22-50-16.237030p00-00_7d12cbea.json'; 'B' = '7d12cbea'; 'C' = 'success'; 'D' = '2025-09-29T22-50-16.237030p00-00'; 'E' = 1; 'F' = 'test'; 'G' = 'openai/gpt-5-mini'; 'H' = 'realpariwise-alltools-gpt5mini-1-test'; 'I' = 'coherence_check, counterfactual_pairs'; 'J' = 0.6666666666666666; 'K' = 0.3333333333333333; 'L' = 0.4444444444444444; 'M' = 0; 'N' = 1; 'O' = 0.3333333333333333; 'P' = 0.5; 'Q' = 3; 'R' = 0.5; 'S' = 0.3333333333333333; 'T' = 0.4; 'U' = 3; 'V' = 0.9166666666666666; 'W' = 0.9777777777777777; 'X' = 0.946236559139785; 'Y' = 45 },
    @{ 'A' = 'run_2025-09-30T00-31-31.473510p00-00_4bfc73d7.json'; 'B' = '4bfc73d7'; 'C' = 'success'; 'D' = '2025-09-30T00-31-31.473510p00-00'; 'E' = 50; 'F' = 'test'; 'G' = 'openai/gpt-5-mini'; 'H' = 'realpariwise-alltools-gpt5mini-50-test'; 'I' = 'coherence_check, counterfactual_pairs'; 'J' = 0.3872549019607843; 'K' = 0.6171875; 'L' = 0.4759036144578314; 'M' = 1; 'N' = 0.3674418604651163; 'O' = 0.6171875; 'P' = 0.4606413994169096; 'Q' = 128; 'R' = 0.4093264248704663; 'S' = 0.6171875; 'T' = 0.4922118380062304; 'U' = 128; 'V' = 0.9698167445203019; 'W' = 0.9195911413969335; 'X' = 0.9440363763553691; 'Y' = 2935 },
)

$startRow = 31
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $rowNum = $startRow + $i
    $rowData = $newRows[$i]
    $ws.Range("A" + $rowNum).Value = $rowData['A']
    $ws.Range("B" + $rowNum).Value = $rowData['B']
    $ws.Range("C" + $rowNum).Value = $rowData['C']
    $ws.Range("D" + $rowNum).Value = $rowData['D']
    $ws.Range("E" + $rowNum).Value = $rowData['E']
    $ws.Range("F" + $rowNum).Value = $rowData['F']
    $ws.Range("G" + $rowNum).Value = $rowData['G']
    $ws.Range("H" + $rowNum).Value = $rowData['H']
    $ws.Range("I" + $rowNum).Value = $rowData['I']
    $ws.Range("J" + $rowNum).Value = $rowData['J']
    $ws.Range("K" + $rowNum).Value = $rowData['K']
    $ws.Range("L" + $rowNum).Value = $rowData['L']
    $ws.Range("M" + $rowNum).Value = $rowData['M']
    $ws.Range("N" + $rowNum).Value = $rowData['N']
    $ws.Range("O" + $rowNum).Value = $rowData['O']
    $ws.Range("P" + $rowNum).Value = $rowData['P']
    $ws.Range("Q" + $rowNum).Value = $rowData['Q']
    $ws.Range("R" + $rowNum).Value = $rowData['R']
    $ws.Range("S" + $rowNum).Value = $rowData['S']
    $ws.Range("T" + $rowNum).Value = $rowData['T']
    $ws.Range("U" + $rowNum).Value = $rowData['U']
    $ws.Range("V" + $rowNum).Value = $rowData['V']
    $ws.Range("W" + $rowNum).Value = $rowData['W']
    $ws.Range("X" + $rowNum).Value = $rowData['X']
    $ws.Range("Y" + $rowNum).Value = $rowData['Y']
}